$d = $word.ActiveDocument

# Insert two new paragraphs at the very start of the document, before the
# existing first paragraph ("2022.06.29" / "阴"):
#   2022.06.28
#   雨
# Using InsertXML (instead of typing text) lets us reproduce the exact
# run/rPr split that appears in the target OOXML: the leading "2" keeps the
# eastAsia font hint (as it would if typed via an East-Asian IME) while the
# remainder of the date does not, matching the pattern already used by the
# existing "2022.06.29" paragraph further down in the document.
$start = $d.Range(0, 0)
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$xml = "<w:p $wNs><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>2</w:t></w:r><w:r><w:t>022.06.28</w:t></w:r></w:p>" +
       "<w:p $wNs><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>雨</w:t></w:r></w:p>"
$start.InsertXML($xml)
